$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column updates. Several new values (e.g. "3.00", "0.491")
# look numeric, so a plain .Value assignment would have Excel silently
# re-parse them as numbers and drop the trailing zero / exact text.
# Force text storage per cell, then restore the default (unstyled)
# style so no stray number-format is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.452.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.253.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.245.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.491"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.778.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "557.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.552.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.256.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "566.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0455"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0868"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.210.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.283"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0562"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.64"
$ws.Range("D49").Style = "Normal"

# Volume(1h) (E) column updates (plain text, percentages with padding spaces
# so Excel keeps them as text already)
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("E3").Value = "  +6.39%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +4.31%  "
$ws.Range("E6").Value = "  +7.33%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +6.52%  "
$ws.Range("E9").Value = "  +5.16%  "
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("E11").Value = "  +5.39%  "
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("E14").Value = "  +5.21%  "
$ws.Range("E15").Value = "  +6.64%  "
$ws.Range("E16").Value = "  +11.95%  "
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("E18").Value = "  +6.59%  "
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("E20").Value = "  +5.50%  "
$ws.Range("E21").Value = "  +4.65%  "
$ws.Range("E22").Value = "  +7.84%  "
$ws.Range("E23").Value = "  +9.09%  "
$ws.Range("E24").Value = "  +6.14%  "
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +17.83%  "
$ws.Range("E28").Value = "  +8.14%  "
$ws.Range("E29").Value = "  +5.42%  "
$ws.Range("E30").Value = "  +5.91%  "
$ws.Range("E31").Value = "  +5.18%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +5.31%  "
$ws.Range("E34").Value = "  +8.25%  "
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("E36").Value = "  +6.55%  "
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("E38").Value = "  +11.29%  "
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("E40").Value = "  +6.89%  "
$ws.Range("E41").Value = "  +9.94%  "
$ws.Range("E42").Value = "  +10.33%  "
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("E44").Value = "  +13.95%  "
$ws.Range("E45").Value = "  +9.07%  "
$ws.Range("E46").Value = "  +5.03%  "
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  +4.69%  "
$ws.Range("E50").Value = "  +3.10%  "
$ws.Range("E51").Value = "  +7.13%  "
